$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the cleaned dataset.
# "RM 232" is row 26; after it is removed, "SC 92" (originally row 28)
# becomes row 27, so we delete row 27 next.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Re-impute / re-mask individual Column D / Column E values (post row-delete
# row numbers) to match the new error-calculation pass.
$ws.Range("E6").Value = -5.7
$ws.Range("E8").Value = ""
$ws.Range("E12").Value = -5.3
$ws.Range("E14").Value = ""
$ws.Range("E17").Value = -7.3
$ws.Range("E18").Value = -8.5
$ws.Range("E19").Value = ""
$ws.Range("E20").Value = ""
$ws.Range("E23").Value = -7

$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = ""
$ws.Range("D28").Value = ""
$ws.Range("D29").Value = ""
$ws.Range("D30").Value = -13.6
$ws.Range("D32").Value = ""
